$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.370.41"
$ws.Range("E2").Value = "  -2.09%  "

$ws.Range("D3").Value = "1.662.73"
$ws.Range("E3").Value = "  -4.01%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9978"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.09"
$ws.Range("E5").Value = "  -3.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9989"
$ws.Range("E6").Value = "  -0.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4805"
$ws.Range("E7").Value = "  -4.55%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2599"
$ws.Range("E8").Value = "  -4.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06149"
$ws.Range("E9").Value = "  -0.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07070"
$ws.Range("E10").Value = "  -2.29%  "

$ws.Range("D11").Value = "1.651.95"
$ws.Range("E11").Value = "  -4.66%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.67"
$ws.Range("E12").Value = "  -3.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5858"
$ws.Range("E13").Value = "  -10.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.363"
$ws.Range("E14").Value = "  -8.59%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "74.32"
$ws.Range("E15").Value = "  -3.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9989"
$ws.Range("E16").Value = "  -0.22%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9983"
$ws.Range("E17").Value = "  -0.14%  "

$ws.Range("D18").Value = "25.358.58"
$ws.Range("E18").Value = "  -2.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006688"
$ws.Range("E19").Value = "  -2.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.41"
$ws.Range("E20").Value = "  -4.35%  "

$ws.Range("D21").Value = "1.865.76"
$ws.Range("E21").Value = "  -4.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.389"
$ws.Range("E22").Value = "  -4.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.598"
$ws.Range("E23").Value = "  -2.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.306"
$ws.Range("E24").Value = "  -3.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "134.12"
$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.10"
$ws.Range("E26").Value = "  -1.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.390"
$ws.Range("E27").Value = "  -2.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "104.99"
$ws.Range("E28").Value = "  -0.65%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.676"
$ws.Range("E29").Value = "  -6.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.990"

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.609"
$ws.Range("E31").Value = "  -2.45%  "

$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07647"
$ws.Range("E32").Value = "  -6.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04367"
$ws.Range("E33").Value = "  -7.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9982"
$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.596"
$ws.Range("E35").Value = "  -2.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6020"
$ws.Range("E36").Value = "  -1.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9415"
$ws.Range("E37").Value = "  -5.61%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.616"
$ws.Range("E38").Value = "  -4.77%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8501"
$ws.Range("E39").Value = "  -3.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9998"
$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("E41").Value = "  -6.94%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.811"
$ws.Range("E42").Value = "  -7.83%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "98.53"
$ws.Range("E43").Value = "  -3.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3744"
$ws.Range("E44").Value = "  -4.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.685"
$ws.Range("E45").Value = "  -6.80%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1112"
$ws.Range("E46").Value = "  -5.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.191"
$ws.Range("E47").Value = "  -3.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05241"
$ws.Range("E48").Value = "  -0.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.39"
$ws.Range("E49").Value = "  -4.75%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.216"
$ws.Range("E50").Value = "  -2.06%  "

$ws.Range("E51").Value = "  -0.12%  "
